$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# This sheet logs one weekly price observation per row (rows 2..104), each
# data block repeating through a fixed 7-week cycle of Origen/Unidad groups.
# The commit adds this week's newest observation at the top of the log
# (row 13 - right after the most-recent entries already present in rows
# 2-12) and pushes everything that was at row 13 downwards by one row, so
# the oldest row (104) falls off the end into a brand-new row 105.

# Insert a new blank row at 13; this shifts old rows 13..104 down to 14..105
# (and grows the used range / dimension to A1:R105 automatically).
$ws.Rows("13:13").Insert()

# Populate the newly inserted row 13 with this week's new record.
$ws.Range("A13").Value = 4
$ws.Range("B13").Value = 'Feria Lagunitas de Puerto Montt'
$ws.Range("C13").Value = 'Los Lagos'
$ws.Range("D13").Value = 44462
$ws.Range("E13").Value = 10
$ws.Range("F13").Value = 100112009
$ws.Range("G13").Value = 'Acelga'
$ws.Range("H13").Value = 'Sin especificar'
$ws.Range("I13").Value = 'Primera'
$ws.Range("J13").Value = 80
$ws.Range("K13").Value = 4000
$ws.Range("L13").Value = 4000
$ws.Range("M13").Value = 4000
$ws.Range("N13").Value = '$/docena de atados (4 kilos)'
$ws.Range("O13").Value = 'Región del Maule'
$ws.Range("P13").Value = 1000
$ws.Range("Q13").Value = 4
$ws.Range("R13").Value = 'Hortaliza'
